$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33, pushing the existing rows 33:57 down to 34:58
$ws.Rows(33).Insert()

# Populate the new row 33 with the new weekly record (same Mercado/Categoria
# metadata as its neighbours, new date + volume/price figures)
$ws.Range("A33").Value = 1
$ws.Range("B33").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C33").Value = "Arica y Parinacota"
$ws.Range("D33").Value = 45264
$ws.Range("E33").Value = 15
$ws.Range("F33").Value = 100112003
$ws.Range("G33").Value = "Ajo"
$ws.Range("H33").Value = "Chino"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 380
$ws.Range("K33").Value = 18000
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = 19053
$ws.Range("N33").Value = "$/caja 10 kilos"
$ws.Range("O33").Value = "China"
$ws.Range("P33").Value = 1905
$ws.Range("Q33").Value = 10
$ws.Range("R33").Value = "Hortaliza"

# Match the date-formatted style used by the other rows' Fecha column
$ws.Range("D33").NumberFormat = $ws.Range("D34").NumberFormat
